# Generate Report for Handoff
#
# The localization-status report is regenerated: a new handoff bundle GUID
# replaces the old one, the handoff package hash changes, and the
# handoff/handback timestamps advance. This touches:
#   - the three "<guid>.md" / "<guid>.<hash>.<lang>.xlf" source-file-name
#     cells (and their hyperlink display text) on Overview/zh-cn/de-de
#   - the two "Latest Handoff Datetime" cells on zh-cn/de-de
#
# NOTE: the hyperlink *targets* (the underlying URLs, stored in each
# worksheet's .rels) are left exactly as they were - only the visible
# "display" text is refreshed. Because this COM host always appends a new
# hyperlink object instead of editing an existing one's display text in
# place, each sheet's hyperlinks are rebuilt from scratch (delete-all then
# re-add in original order) so the relationship ids/targets come out
# unchanged while only the display strings change.

$wb = $excel.ActiveWorkbook

$oldGuid = "b58df756-4b38-41ee-93f4-2409645db38a"
$newGuid = "bad0d6d0-02de-459b-b110-cbce4f192c5f"

$oldHash = "06a61feb744cbc0ec2b9a392867bb2efe271309b"
$newHash = "416c687a14111486b77eae7aca30ef6f491b86c4"

$oldZhDate = "2016-03-08 23:22:22"
$newZhDate = "2016-03-08 23:23:13"

$oldDeDate = "2016-03-08 23:22:30"
$newDeDate = "2016-03-08 23:23:22"

$newMd = "$newGuid.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4e115c2b831e9b8ef3e6ccc4f7cf1522de7aa59a/e2e/$oldGuid.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4e115c2b831e9b8ef3e6ccc4f7cf1522de7aa59a/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1bb4e564e92c73d22079ca579addbcce8aa9cd01/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca0dc2686d6c8a1fddd692ac3b84cd078f91e7be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config")

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", ".localization-config")

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", ".localization-config")
